$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values would be auto-detected as numbers by Excel; force Text
# number format first so COM stores them as text, matching the source format.
$forceTextCells = @("D5","D7","D8","D9","D10","D12","D14","D15","D16","D18","D21","D22","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D37","D38","D39","D40","D41","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated price / volume figures cell by cell.
$ws.Range("D2").Value = "27.434.15"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "1.861.02"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "311.27"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4773"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").Value = "0.07319"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "0.9311"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  +5.27%  "
$ws.Range("D12").Value = "0.07792"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "1.856.12"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "5.442"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").Value = "6.551"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "90.12"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "0.000008821"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "27.476.74"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "14.62"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "5.097"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "155.27"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "18.47"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "2.010"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "115.48"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").Value = "4.951"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "0.08895"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "3.332"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").Value = "1.206"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").Value = "0.7533"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Value = "4.598"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "0.02040"
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("D38").Value = "0.5557"
$ws.Range("E38").Value = "  +5.83%  "
$ws.Range("D39").Value = "0.05272"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "2.990"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "7.038"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").Value = "0.1522"
$ws.Range("D44").Value = "0.4879"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").Value = "10.63"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "1.669"
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("D48").Value = "103.06"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "67.42"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "0.06091"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "0.9121"
$ws.Range("E51").Value = "  +2.56%  "
